$d = $word.ActiveDocument

# Find the run containing the target sentence
$target = "Tests should take around 15 mins to run on a medium performance desktop. "

$found = $d.Content.Find.Execute("Tests should take around 15 mins to run on a medium performance desktop. ", $true, $false, $false, $false, $false, $true, 1, $false, "Tests should take around 15-30 mins to run on a medium performance desktop. ", 2)

Write-Host "Find result: $found"
